$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.378.67'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '1.884.58'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.695'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '246.45'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.71%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.23'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.357'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0750'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '13.52'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.777'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +8.60%  '
$ws.Range('D14').Value = '2.161.78'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '1.918.39'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D17').Value = '35.387.26'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '73.87'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = '0.0₃0829'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.84'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.24'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.62'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +10.18%  '
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  -3.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.02'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.67'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.39%  '
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.20'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.86'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.48'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.29%  '
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0739'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +11.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.38'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('E40').Value = '  +4.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.60'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('E42').Value = '  -0.57%  '
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('D44').Value = '1.312.05'
$ws.Range('E44').Value = '  +1.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0807'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +6.23%  '
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.07'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('E49').Value = '  -1.63%  '
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').Value = '2.068.13'
$ws.Range('E51').Value = '  +0.15%  '
